$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$teams = @("POR","CLE","DAL","MIA","OKC","ATL","WAS","MIL","LAC","SAS","DET","ORL","UTA","MEM","HOU","NOP","DEN","LAL","GSW","IND","CHO","CHI","PHI","BOS","BRK","TOR","SAC","PHO","NYK","MIN")
$pers = @(14.39230769230769,12.59230769230769,15.30833333333333,11.93125,15.66363636363636,15.21428571428571,11.34,13.35384615384615,12.45882352941176,14.98666666666667,13.64,12.16,12.56666666666667,15.1,12.78333333333333,12.67333333333333,12.30833333333334,11.45882352941176,15.46,14.14666666666667,13.06428571428571,12.81428571428571,8.507142857142858,13.08,12.25333333333333,12.97333333333333,16.7625,9.936363636363636,12.95,13.10625)

for ($i = 0; $i -lt $teams.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $teams[$i]
    $ws.Cells.Item($row, 3).Value = $pers[$i]
}
